$d = $word.ActiveDocument

# The "_GoBack" bookmark (Word's "last edit location" marker) was sitting
# around the "<Proposal Description>" run from a previous edit. Remove it
# from there -- it needs to move to the signature block below, where the
# actual edit (fixing the <Casework Officer> run's font) now happens.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Find the "Key_body"-styled paragraph that holds the "<Casework Officer>"
# placeholder used as the letter's signature line (the one right under
# "Yours sincerely" -- distinct from the "Key_HeadDetails" one near the
# top of the letter that also mentions "Casework Officer").
$target = $null
foreach ($p in $d.Paragraphs) {
    if (($p.Range.Text -like "*<Casework Officer>*") -and `
        ($p.Range.Style.NameLocal -eq "Key_body")) {
        $target = $p
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # The signature format was wrong: it inherited "Source Sans Pro Light"
    # from the Key_body style. Force the correct "Source Sans Pro" on both
    # the run and the paragraph mark (matches Word's own behaviour when
    # applying a font to a paragraph range that includes the pilcrow).
    $r.Font.Name = "Source Sans Pro"

    # Re-create "_GoBack" around the "<Casework Officer>" text, marking it
    # as the most recently edited spot in the document.
    $txt = $r.Text
    $markerLen = $txt.IndexOf("<Casework Officer>") + "<Casework Officer>".Length
    $bmRange = $d.Range($r.Start, $r.Start + $markerLen)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
